$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timelog entries continuing "thesis chapter previous work" (2015-04-29 / 2015-04-30),
# finishing that chapter of work. Serial dates 42123 = 2015-04-29, 42124 = 2015-04-30.

$rows = @(
    @{ Row = 107; Date = 42123; From = 0.29166666666666669; To = 0.52083333333333337 },
    @{ Row = 108; Date = 42123; From = 0.5625;              To = 0.72916666666666663 },
    @{ Row = 109; Date = 42124; From = 0.47916666666666669; To = 0.53125 },
    @{ Row = 110; Date = 42124; From = 0.57291666666666663; To = 0.75 }
)

# Reuse the existing date / time number formats (copy format-only from a
# fully populated row) so no duplicate style entries get created.
$xlPasteFormats = -4122

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A104").Copy()
    $ws.Range("A$row").PasteSpecial($xlPasteFormats)
    $ws.Range("A$row").Value = $r.Date

    $ws.Range("B104").Copy()
    $ws.Range("B$row").PasteSpecial($xlPasteFormats)
    $ws.Range("B$row").Value = $r.From

    $ws.Range("C104").Copy()
    $ws.Range("C$row").PasteSpecial($xlPasteFormats)
    $ws.Range("C$row").Value = $r.To

    $ws.Range("E$row").Value = "thesis chapter previous work"
}

$excel.CutCopyMode = $false

$wb.Application.Calculate()

$ws.Range("E114").Select()
